$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New event rows (kosice, romania, flensburg, almahatta sitta) computed
# empirically from the CNEOS fireball data, appended below the existing
# table (rows 2-15).

$ws.Range("A16").Value = "kosice"
$ws.Range("B16").Value = 153000000000
$ws.Range("C16").Value = 0.44
$ws.Range("D16").Formula = "=B16/ (4184000000000 * C16)"

$ws.Range("A17").Value = "romania"
$ws.Range("B17").Value = 136000000000
$ws.Range("C17").Value = 0.4
$ws.Range("D17").Formula = "=B17/ (4184000000000 * C17)"

$ws.Range("A18").Value = "flensburg"
$ws.Range("B18").Value = 169000000000
$ws.Range("C18").Value = 0.48
$ws.Range("D18").Formula = "=B18/ (4184000000000 * C18)"

$ws.Range("A19").Value = "almahatta sitta"
$ws.Range("B19").Value = 395000000000
$ws.Range("C19").Value = 1
$ws.Range("D19").Formula = "=B19/ (4184000000000 * C19)"

# Match the new cell style introduced alongside the data (wrapped text on
# the radiated(J)/total(kT) columns for the freshly added rows).
$ws.Range("B16:C19").WrapText = $true
